$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per the new ANFIS/RNN/ARIMA comparison data layout
$ws.Range("A1").Value = 26216.599999999991
$ws.Range("B1").Value = 23176.999999999996
$ws.Range("C1").Value = 22943.100000000002
$ws.Range("A2").Value = 6668
$ws.Range("B2").Value = 2856
$ws.Range("C2").Value = 2973
$ws.Range("A3").Value = 2104
$ws.Range("B3").Value = 2515
$ws.Range("C3").Value = 2611
$ws.Range("A4").Value = 3637
$ws.Range("B4").Value = 5012
$ws.Range("C4").Value = 5175
$ws.Range("A5").Value = 886.08300000000008
$ws.Range("B5").Value = 1925.7290000000005
$ws.Range("C5").Value = 2129.4700000000003
$ws.Range("A6").Value = 6659.7100000000019
$ws.Range("B6").Value = 3682.746000000001
$ws.Range("C6").Value = 4209.5439999999999
$ws.Range("A7").Value = 24469.299999999992
$ws.Range("B7").Value = 21674.298999999999
$ws.Range("C7").Value = 21170.699999999997
$ws.Range("A8").Value = 5513
$ws.Range("B8").Value = 2471
$ws.Range("C8").Value = 2556
$ws.Range("A9").Value = 3489
$ws.Range("B9").Value = 2788
$ws.Range("C9").Value = 2807
$ws.Range("A10").Value = 3298
$ws.Range("B10").Value = 4467
$ws.Range("C10").Value = 4646
$ws.Range("A11").Value = 769.35000000000014
$ws.Range("B11").Value = 2061.8949999999991
$ws.Range("C11").Value = 2227.2549999999992
$ws.Range("A12").Value = 6190.8819999999969
$ws.Range("B12").Value = 3433.9150000000004
$ws.Range("C12").Value = 4066.0949999999993
$ws.Range("A13").Value = 24606.971999999965
$ws.Range("B13").Value = 20376.971999999998
$ws.Range("C13").Value = 20264.103999999999
$ws.Range("A14").Value = 5370
$ws.Range("B14").Value = 2316
$ws.Range("C14").Value = 2309
$ws.Range("A15").Value = 3586
$ws.Range("B15").Value = 2166
$ws.Range("C15").Value = 1995
$ws.Range("A16").Value = 2963
$ws.Range("B16").Value = 3986
$ws.Range("C16").Value = 3946
$ws.Range("A17").Value = 1680.8300000000006
$ws.Range("B17").Value = 2365.848
$ws.Range("C17").Value = 2592.3199999999993
$ws.Range("A18").Value = 6418.5089999999991
$ws.Range("B18").Value = 3538.2750000000015
$ws.Range("C18").Value = 4137.1259999999975
$ws.Range("A19").Value = 25596.69999999999
$ws.Range("B19").Value = 14632.499999999998
$ws.Range("C19").Value = 14228.300000000005
$ws.Range("A20").Value = 4364
$ws.Range("B20").Value = 1440
$ws.Range("C20").Value = 1579
$ws.Range("A21").Value = 1373.362499999999
$ws.Range("B21").Value = 974.58799999999962
$ws.Range("C21").Value = 904.50499999999965
$ws.Range("A22").Value = 2781
$ws.Range("B22").Value = 2733
$ws.Range("C22").Value = 2582
$ws.Range("A23").Value = 2226.7299999999991
$ws.Range("B23").Value = 1618.7199999999998
$ws.Range("C23").Value = 1595.3999999999992
$ws.Range("A24").Value = 6446.7100000000028
$ws.Range("B24").Value = 2889.52
$ws.Range("C24").Value = 3212.8600000000006
$ws.Range("A25").Value = 25326.990247422673
$ws.Range("B25").Value = 18115.577752577326
$ws.Range("C25").Value = 17133.401000000002
$ws.Range("A26").Value = 5873
$ws.Range("B26").Value = 2606
$ws.Range("C26").Value = 2656
$ws.Range("A27").Value = 1437.6405000000002
$ws.Range("B27").Value = 1305.6855000000005
$ws.Range("C27").Value = 1296.9269999999997
$ws.Range("A28").Value = 2865
$ws.Range("B28").Value = 5073
$ws.Range("C28").Value = 4456
$ws.Range("A29").Value = 2168.4419999999991
$ws.Range("B29").Value = 1919.04
$ws.Range("C29").Value = 1950.3590000000004
$ws.Range("A30").Value = 5620.2959999999975
$ws.Range("B30").Value = 3200.3049999999994
$ws.Range("C30").Value = 3448.8949999999995
$ws.Range("A31").Value = 25868.700000000001
$ws.Range("B31").Value = 20778.599999999991
$ws.Range("C31").Value = 19379.200000000012
$ws.Range("A32").Value = 6580
$ws.Range("B32").Value = 3214
$ws.Range("C32").Value = 3375
$ws.Range("A33").Value = 1512.3250000000003
$ws.Range("B33").Value = 1711.5965000000003
$ws.Range("C33").Value = 1658.5349999999992
$ws.Range("A34").Value = 3529
$ws.Range("B34").Value = 6042
$ws.Range("C34").Value = 5353
$ws.Range("A35").Value = 1903.3499999999999
$ws.Range("B35").Value = 1415.1059999999993
$ws.Range("C35").Value = 1425.0379999999998
$ws.Range("A36").Value = 5041.8729999999978
$ws.Range("B36").Value = 3122.7190000000005
$ws.Range("C36").Value = 3308.4050000000007
$ws.Range("A37").Value = 31211.599999999955
$ws.Range("B37").Value = 22784.5
$ws.Range("C37").Value = 20634.899999999991
$ws.Range("A38").Value = 8977
$ws.Range("B38").Value = 3985
$ws.Range("C38").Value = 4141
$ws.Range("A39").Value = 5059.6164999999946
$ws.Range("B39").Value = 3905.7699999999977
$ws.Range("C39").Value = 3782.1229999999991
$ws.Range("A40").Value = 6625
$ws.Range("B40").Value = 7274
$ws.Range("C40").Value = 6539
$ws.Range("A41").Value = 1969.172
$ws.Range("B41").Value = 892.43800000000022
$ws.Range("C41").Value = 945.06500000000028
$ws.Range("A42").Value = 5868.4619999999932
$ws.Range("B42").Value = 3155.5130000000022
$ws.Range("C42").Value = 3470.2359999999999
$ws.Range("A43").Value = 25209.899999999987
$ws.Range("B43").Value = 20057.599999999999
$ws.Range("C43").Value = 18162.5
$ws.Range("A44").Value = 8247
$ws.Range("B44").Value = 4344
$ws.Range("C44").Value = 4509
$ws.Range("A45").Value = 5346.5099999999984
$ws.Range("B45").Value = 4057.3584999999989
$ws.Range("C45").Value = 4199.5489999999991
$ws.Range("A46").Value = 6528
$ws.Range("B46").Value = 7836
$ws.Range("C46").Value = 6956
$ws.Range("A47").Value = 1541.1400000000006
$ws.Range("B47").Value = 774.31000000000063
$ws.Range("C47").Value = 781.94999999999982
$ws.Range("A48").Value = 5499.7400000000034
$ws.Range("B48").Value = 3429.6760000000004
$ws.Range("C48").Value = 3680.5840000000017
$ws.Range("A49").Value = 39587.899999999943
$ws.Range("B49").Value = 14817.300000000005
$ws.Range("C49").Value = 14851.200000000012
$ws.Range("A50").Value = 7327
$ws.Range("B50").Value = 3234
$ws.Range("C50").Value = 3381
$ws.Range("A51").Value = 3887.4975000000031
$ws.Range("B51").Value = 2863.0859999999971
$ws.Range("C51").Value = 3005.5939999999991
$ws.Range("A52").Value = 5609
$ws.Range("B52").Value = 5904
$ws.Range("C52").Value = 5267
$ws.Range("A53").Value = 2235.9100000000017
$ws.Range("B53").Value = 1762.1320000000005
$ws.Range("C53").Value = 1902.9730000000004
$ws.Range("A54").Value = 4869.9800000000023
$ws.Range("B54").Value = 2704.0450000000005
$ws.Range("C54").Value = 2941.7550000000001
$ws.Range("A55").Value = 37817.099999999969
$ws.Range("B55").Value = 14558.898999999996
$ws.Range("C55").Value = 14776.000000000007
$ws.Range("A56").Value = 5831
$ws.Range("B56").Value = 2479
$ws.Range("C56").Value = 2583
$ws.Range("A57").Value = 3381.9580000000005
$ws.Range("B57").Value = 2380.5339999999983
$ws.Range("C57").Value = 2428.5495000000019
$ws.Range("A58").Value = 4039
$ws.Range("B58").Value = 4239
$ws.Range("C58").Value = 3796
$ws.Range("A59").Value = 2284.5360000000001
$ws.Range("B59").Value = 1891.0830000000001
$ws.Range("C59").Value = 2068.4569999999999
$ws.Range("A60").Value = 4328.6750000000002
$ws.Range("B60").Value = 2233.5149999999994
$ws.Range("C60").Value = 2605.3999999999987
$ws.Range("A61").Value = 0
$ws.Range("B61").Value = 0
$ws.Range("C61").Value = 0
$ws.Range("A62").Value = 0
$ws.Range("B62").Value = 0
$ws.Range("C62").Value = 0
$ws.Range("A63").Value = 0
$ws.Range("B63").Value = 0
$ws.Range("C63").Value = 0
$ws.Range("A64").Value = 0
$ws.Range("B64").Value = 0
$ws.Range("C64").Value = 0
$ws.Range("A65").Value = 0
$ws.Range("B65").Value = 0
$ws.Range("C65").Value = 0
$ws.Range("A66").Value = 0
$ws.Range("B66").Value = 0
$ws.Range("C66").Value = 0
$ws.Range("A67").Value = 0
$ws.Range("B67").Value = 0
$ws.Range("C67").Value = 0
$ws.Range("A68").Value = 0
$ws.Range("B68").Value = 0
$ws.Range("C68").Value = 0
$ws.Range("A69").Value = 0
$ws.Range("B69").Value = 0
$ws.Range("C69").Value = 0
$ws.Range("A70").Value = 0
$ws.Range("B70").Value = 0
$ws.Range("C70").Value = 0
$ws.Range("A71").Value = 0
$ws.Range("B71").Value = 0
$ws.Range("C71").Value = 0
$ws.Range("A72").Value = 0
$ws.Range("B72").Value = 0
$ws.Range("C72").Value = 0

# Update column widths
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(2).ColumnWidth = 10.666666666666666
$ws.Columns.Item(3).ColumnWidth = 8.666666666666666
